$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(326, 1).Value = 325
$ws.Cells.Item(326, 2).Value = 'キーワード'
$ws.Cells.Item(326, 3).Value = '【指定攻撃】を持つフォロワーは、相手の【守護】を持つアクト状態のフォロワーを無視して相手の他のフォロワーに攻撃できますか？'
$ws.Cells.Item(326, 4).Value = 'いいえ。必ずアクト状態の【守護】を持つフォロワーを選ぶ必要があります。'

$ws.Cells.Item(327, 1).Value = 326
$ws.Cells.Item(327, 2).Value = 'アーチャー'
$ws.Cells.Item(327, 3).Value = '『エンシェントエルフ』の《ファンファーレ》能力のコストで『アーチャー』を選択した場合、『アーチャー』の能力はプレイしますか？'
$ws.Cells.Item(327, 4).Value = 'はい。一度誘発した自動能力は、そのカードが場を離れても能力をプレイできます。'

$ws.Cells.Item(328, 1).Value = 327
$ws.Cells.Item(328, 2).Value = 'クラフトウォーロック'
$ws.Cells.Item(328, 3).Value = '『クラフトウォーロック(EVOLVE)』の能力で「【スタック】+1する」は自分の【スタック】を持つアミュレット全ての【スタック】を+できますか？'
$ws.Cells.Item(328, 4).Value = 'いいえ。【スタック】を持つアミュレット1枚を選択し、そのカードのスタックカウンターを+1します。
'

$ws.Cells.Item(329, 1).Value = 328
$ws.Cells.Item(329, 2).Value = 'キーワード'
$ws.Cells.Item(329, 3).Value = '【スタック】を持つアミュレットを『エクスキューション』で破壊した場合、【スタック】能力でスタックカウンターを取り除くことで、場に残せますか？
'
$ws.Cells.Item(329, 4).Value = 'はい。自分や相手のカード問わず、場を離れる時に【スタック】能力でスタックカウンターを取り除くことで場に残せます。
'

$ws.Cells.Item(330, 1).Value = 329
$ws.Cells.Item(330, 2).Value = '次元の魔女・ドロシー'
$ws.Cells.Item(330, 3).Value = '『次元の魔女・ドロシー(EVOLVE)』の【進化時】能力でコストが-5された『次元の超越』をプレイする際、さらに『次元の超越』の能力でコストを7にする場合、コストは7から-5をし、コスト2でプレイできますか？
'
$ws.Cells.Item(330, 4).Value = 'はい。コスト2でプレイすることができます。
'

$ws.Cells.Item(331, 1).Value = 330
$ws.Cells.Item(331, 2).Value = 'キーワード'
$ws.Cells.Item(331, 3).Value = '【進化時】能力が複数あるカードに進化した場合、その能力は全てプレイしますか？
'
$ws.Cells.Item(331, 4).Value = 'はい。【進化時】能力が複数ある場合でも、全てプレイします。
'

$ws.Cells.Item(332, 1).Value = 331
$ws.Cells.Item(332, 2).Value = '裁きの悪魔'
$ws.Cells.Item(332, 3).Value = 'このフォロワーが相手のフォロワーと交戦した際に、お互いのフォロワーが破壊された場合、このフォロワーの能力はプレイできますか？
'
$ws.Cells.Item(332, 4).Value = 'はい。交戦でお互いが破壊された場合でも「破壊されたとき」を満たしているため、プレイすることができます。
'

$ws.Cells.Item(333, 1).Value = 332
$ws.Cells.Item(333, 2).Value = '享楽の悪魔'
$ws.Cells.Item(333, 3).Value = 'このフォロワーの能力は、【真紅】状態でない場合でも手札を1枚捨てますか？
'
$ws.Cells.Item(333, 4).Value = 'いいえ。「【真紅】状態なら」の以降の能力は、条件を満たさない場合、プレイできません。
'

$ws.Cells.Item(334, 1).Value = 333
$ws.Cells.Item(334, 2).Value = 'エイラの祈祷'
$ws.Cells.Item(334, 3).Value = 'このアミュレットの能力は相手のターンでもプレイできますか？
'
$ws.Cells.Item(334, 4).Value = 'はい。相手のターンでも条件を満たせば、能力はプレイできます。
'

$ws.Cells.Item(335, 1).Value = 334
$ws.Cells.Item(335, 2).Value = '神域の守護者'
$ws.Cells.Item(335, 3).Value = 'このアミュレットの能力は相手のターンでもプレイできますか？
'
$ws.Cells.Item(335, 4).Value = 'いいえ。このアミュレットの能力は自分のターンのみプレイできます。
'

$ws.Cells.Item(336, 1).Value = 335
$ws.Cells.Item(336, 2).Value = 'エイラの祈祷'
$ws.Cells.Item(336, 3).Value = '『ダークオファリング』の能力で『エイラの祈祷』を選択し、そのカードを破壊して、リーダーの《体力》を+した場合、『エイラの祈祷』の能力はプレイできますか？
'
$ws.Cells.Item(336, 4).Value = 'いいえ。『ダークオファリング』で自分のリーダーを《体力》+3したとき、『エイラの祈祷』が場にないため「リーダーの《体力》を+したとき」の能力はプレイできません。
'

$ws.Cells.Item(337, 1).Value = 336
$ws.Cells.Item(337, 2).Value = 'ダークエンジェル・オリヴィエ'
$ws.Cells.Item(337, 3).Value = '先攻のプレイヤーが『ダークエンジェル・オリヴィエ』の能力でEPを+1できますか？
'
$ws.Cells.Item(337, 4).Value = 'はい。先攻後攻問わず、全てのプレイヤーは能力でEPを得ます。
'

$ws.Cells.Item(338, 1).Value = 337
$ws.Cells.Item(338, 2).Value = 'ダークエンジェル・オリヴィエ'
$ws.Cells.Item(338, 3).Value = 'EPが3ある状態で『ダークエンジェル・オリヴィエ』の能力でEPを+1できますか？
'
$ws.Cells.Item(338, 4).Value = 'はい。EPは上限がないため、EPが3ある場合も新たにEPを得ます。
'
